# Lecture partielle de l'EDT M1 MIAGE.
#
# The sheet is a weekly schedule export. Column A holds the date of each
# day-block (numeric, date-formatted) and column B holds the French name
# of that day of week (text, e.g. "jeudi"). This shifts the whole
# schedule forward by 1096 days (= 3 years, e.g. 2023 -> 2026) and fixes
# up the day-of-week labels in column B to match the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayOffset = 1096

# French day-of-week names, Excel's VBA Weekday() convention: 1=dimanche .. 7=samedi
$dayNames = @("dimanche", "lundi", "mardi", "mercredi", "jeudi", "vendredi", "samedi")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2

    if ($aVal -ne $null -and ($aVal -is [double] -or $aVal -is [int])) {
        $newDate = $aVal + $dayOffset
        $aCell.Value = $newDate

        # weekday index: days since 1899-12-30 (serial 0) mod 7, serial 1 = Monday 1900-01-01
        # Excel serial 2 (1900-01-01) was a Monday -> weekday 2 in the 1-based (dimanche=1) scheme
        $weekdayIndex = (([int]$newDate + 6) % 7)
        $dayName = $dayNames[$weekdayIndex]

        $bCell = $ws.Cells.Item($r, 2)
        if ($bCell.Value2 -ne $null) {
            $bCell.Value = $dayName
        }
    }
}
